$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.056.80"
$ws.Range("D3").Value = "1.747.57"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'226.04"
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("D6").Value = "'0.5808"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.2711"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "'23.18"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.06593"
$ws.Range("E10").Value = "  -4.90%  "
$ws.Range("D11").Value = "'0.07508"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").Value = "1.748.63"
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("D13").Value = "'4.712"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "'0.6030"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "1.985.28"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("D16").Value = "'73.93"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").Value = "'0.000008626"
$ws.Range("E17").Value = "  -10.91%  "
$ws.Range("D18").Value = "28.046.86"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").Value = "'5.319"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'204.92"
$ws.Range("E21").Value = "  -4.79%  "
$ws.Range("D22").Value = "'11.27"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").Value = "'6.652"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'150.44"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("D26").Value = "'8.013"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "'0.1236"
$ws.Range("E27").Value = "  -3.44%  "
$ws.Range("D28").Value = "'16.08"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("D29").Value = "'1.388"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "'0.06099"
$ws.Range("E30").Value = "  -4.50%  "
$ws.Range("D31").Value = "'1.385"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").Value = "'3.729"
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("D33").Value = "'3.710"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "'1.678"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("E35").Value = "  -4.41%  "
$ws.Range("D36").Value = "'0.6351"
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("D37").Value = "'2.433"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("D38").Value = "'2.663"
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").Value = "'0.01678"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").Value = "'6.276"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").Value = "1.126.60"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("D42").Value = "'0.8644"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "'99.62"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "1.898.53"
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("D46").Value = "'59.10"
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("D47").Value = "'1.578"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "'0.00000000108"
$ws.Range("E48").Value = "  -4.43%  "
$ws.Range("D49").Value = "'8.275"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "'0.05397"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("B51").Value = "Frax"
$ws.Range("C51").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D51").Value = "'1.010"
$ws.Range("E51").Value = "  +0.46%  "

# Reset style on cells that were forced to text via apostrophe prefix,
# to avoid introducing a quotePrefix style not present in the original file.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
